$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 137 (shifts existing rows 137:149 down to 138:150)
$ws.Rows(137).Insert()

# Populate the newly inserted row with the new weekly Mango price record
$ws.Range("A137").Value = 5
$ws.Range("B137").Value = "Macroferia Regional de Talca"
$ws.Range("C137").Value = "Maule"
$ws.Range("D137").Value = 44826
$ws.Range("E137").Value = 7
$ws.Range("F137").Value = "Fruta"
$ws.Range("G137").Value = 100108
$ws.Range("H137").Value = "Tropicales y subtropicales"
$ws.Range("I137").Value = 100108002
$ws.Range("J137").Value = "Mango"
$ws.Range("K137").Value = "Sin especificar"
$ws.Range("L137").Value = "Primera"
$ws.Range("M137").Value = 228
$ws.Range("N137").Value = 9000
$ws.Range("O137").Value = 9000
$ws.Range("P137").Value = 9000
$ws.Range("Q137").Value = "$/bandeja 4 kilos"
$ws.Range("R137").Value = "Brasil"
$ws.Range("S137").Value = 2250
$ws.Range("T137").Value = 4
